$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append test flag info to the existing strings
$ws.Range("E2").Value = "environment=0-100,Init=0-0;test=true"
$ws.Range("H2").Value = "Init+1;Time=0; test=true; apa=false"

# Move the active selection
$ws.Range("E2").Select()
